# Recipes_addimage_url.xlsx — add a "Resized Image URL" column (K) that
# mirrors the existing "Image URL" column (J): for every row whose J cell
# holds an image URL like ".../<n>.png", K gets ".../<n>_resized.png".
# Rows whose J cell is blank get a blank K cell too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 147

# Seed column K with J's formatting (header style, borders, etc.) for the
# whole used range in one shot; values get overwritten below.
$ws.Range("J1:J147").Copy($ws.Range("K1:K147"))

$headerCell = $ws.Cells.Item(1, 11)
$headerCell.Value2 = "Resized Image URL"

for ($row = 2; $row -le $lastRow; $row++) {
    $srcCell = $ws.Cells.Item($row, 10)   # column J
    $dstCell = $ws.Cells.Item($row, 11)   # column K
    $url = $srcCell.Value2

    if ($url -and $url.Length -gt 0) {
        $slashIdx = $url.LastIndexOf("/")
        $dotIdx = $url.LastIndexOf(".")
        $baseName = $url.Substring($slashIdx + 1, $dotIdx - $slashIdx - 1)
        $ext = $url.Substring($dotIdx)
        $prefix = $url.Substring(0, $slashIdx + 1)
        $dstCell.Value2 = $prefix + $baseName + "_resized" + $ext
    }
}
